$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.317.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.72%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.582.02'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.37%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.33%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.89'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.73%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.591'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.594.87'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.68%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.68'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.65%  '

$ws.Range("E11").Value = '  +0.78%  '

$ws.Range("E12").Value = '  +5.62%  '

$ws.Range("E13").Value = '  +4.72%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.039.83'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.60%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '59.322.24'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.60%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.01'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.82%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000137'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.584.83'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.54'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.98%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '339.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.38%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.49%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.22%  '

$ws.Range("E23").Value = '  -0.33%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.483'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +8.60%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.01%  '

$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("E27").Value = '  -1.55%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.45'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0766'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.45%  '

$ws.Range("E30").Value = '  -0.02%  '

$ws.Range("E31").Value = '  +0.79%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '157.92'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.29'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.81%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.913'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.30%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.15'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.77%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '37.51'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.09%  '

$ws.Range("E39").Value = '  +0.61%  '

$ws.Range("E40").Value = '  -4.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.66'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '291.05'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '136.28'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.82%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0970'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.09%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.598'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.65'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.49%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0533'
$ws.Range("D48").Style = "Normal"

$ws.Range("E49").Value = '  +1.93%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.958.81'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.84%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.60'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.92%  '

